# UML Project check list.xlsx - "added activity diagram + updated project plan"
#
# Marks the "x" (done) column for the ACTIVITY DIAGRAM checklist rows that
# were completed (rows 48, 49, 52, 53, 54 and 57 in Sheet1, column C),
# mirroring the pattern already used throughout the rest of the sheet.
# Also moves the saved selection/scroll position to where the user ended
# up working (cell C48).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$doneRows = @(48, 49, 52, 53, 54, 57)
foreach ($r in $doneRows) {
    $ws.Cells.Item($r, 3).Value = "x"
}

# Reflect the new working position: selection on C48 (Excel also records
# this as the sheet's saved scroll/active cell).
$ws.Range("C48").Select()
